$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '80.952.28'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +2.60%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.138.75'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.67%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '207.89'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.16%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '616.47'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.36%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.280'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +23.72%  '
$ws.Range("E8").Value = '  -0.05%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.575'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '3.132.14'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.91%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.571'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.44%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0000250'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +11.78%  '
$ws.Range("E13").Value = '  -0.22%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.25'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.22%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '3.713.37'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.73%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '31.23'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.66%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '80.767.20'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.45%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.137.75'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.62%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '3.12'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +9.35%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.80'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.50%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '428.49'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.15%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '8.92'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.89%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.05'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +3.02%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '7.15'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +4.58%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '5.16'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +9.03%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.292.83'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.74%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '75.47'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.23%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.74'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.38%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.40%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0000120'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +5.19%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.14%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '8.89'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.80%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '557.52'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +8.39%  '
$ws.Range("E34").Value = '  -0.23%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.152'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +13.71%  '
$ws.Range("E36").Value = '  +7.54%  '
$ws.Range("E37").Value = '  -1.28%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '22.47'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("E39").Value = '  -0.12%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '6.03'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +11.67%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.403'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.50%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '20.72'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.63%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.01'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +13.33%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +21.62%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '159.00'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("E46").Value = '  -0.01%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '186.59'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.97%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '44.47'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +4.20%  '
$ws.Range("E49").Value = '  +1.48%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.759'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -5.75%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '25.62'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +3.46%  '
